# Deploy the implementation guide:
# - bump the ValueSet "Status" from "active" to "draft"
# - bump the ValueSet "Date" to the new publication timestamp
#
# These two fields live on the "Metadata" sheet of the generated FHIR
# ValueSet spreadsheet, in column B next to their "Status"/"Date" labels
# in column A (rows 6 and 8 respectively).

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Metadata") {
        $ws = $sheet
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: 2023-05-12T12:33:13+00:00 -> 2023-08-01T16:12:28+00:00
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"
